$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.961.47'
$ws.Range("E2").Value = '  +0.07%  '

$ws.Range("D3").Value = '2.792.11'
$ws.Range("E3").Value = '  -1.51%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '358.72'
$ws.Range("E5").Value = '  +0.54%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.89'
$ws.Range("E6").Value = '  -2.39%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.561'
$ws.Range("E7").Value = '  -0.93%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.592'
$ws.Range("E9").Value = '  -1.79%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.38'
$ws.Range("E10").Value = '  -1.91%  '

$ws.Range("E11").Value = '  +1.98%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0851'
$ws.Range("E12").Value = '  -0.77%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.56'
$ws.Range("E13").Value = '  -2.54%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.61'
$ws.Range("E14").Value = '  -2.07%  '

$ws.Range("D15").Value = '3.227.91'
$ws.Range("E15").Value = '  -1.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.949'
$ws.Range("E16").Value = '  +1.80%  '

$ws.Range("D17").Value = '2.768.35'
$ws.Range("E17").Value = '  -1.79%  '

$ws.Range("D18").Value = '51.894.67'
$ws.Range("E18").Value = '  +0.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.46'
$ws.Range("E19").Value = '  -1.04%  '

$ws.Range("E20").Value = '  -2.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.25'
$ws.Range("E21").Value = '  -1.49%  '

$ws.Range("D22").Value = '0.0₃0978'
$ws.Range("E22").Value = '  -1.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '270.79'
$ws.Range("E23").Value = '  +0.46%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.15'
$ws.Range("E24").Value = '  +0.12%  '

$ws.Range("E25").Value = '  -1.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.49'
$ws.Range("E26").Value = '  -2.13%  '

$ws.Range("E27").Value = '  +0.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.165'
$ws.Range("E28").Value = '  +17.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.26'
$ws.Range("E29").Value = '  -0.59%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.17'
$ws.Range("E30").Value = '  -4.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '52.16'
$ws.Range("E31").Value = '  -0.82%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '34.94'
$ws.Range("E32").Value = '  -1.61%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0465'
$ws.Range("E33").Value = '  -2.31%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.99'
$ws.Range("E34").Value = '  +1.38%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0849'
$ws.Range("E35").Value = '  +0.49%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.19'
$ws.Range("E36").Value = '  -4.40%  '

$ws.Range("E37").Value = '  +0.02%  '

$ws.Range("E38").Value = '  +0.87%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.21'
$ws.Range("E39").Value = '  -2.37%  '

$ws.Range("E40").Value = '  -3.57%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.60'
$ws.Range("E41").Value = '  +1.86%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.115'
$ws.Range("E42").Value = '  -1.78%  '

$ws.Range("E43").Value = '  -2.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.42'
$ws.Range("E44").Value = '  -3.73%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.84'
$ws.Range("E45").Value = '  -6.03%  '

$ws.Range("D46").Value = '2.078.45'
$ws.Range("E46").Value = '  -1.09%  '

$ws.Range("E47").Value = '  -2.57%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.25'
$ws.Range("E48").Value = '  -0.41%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.76'
$ws.Range("E49").Value = '  -3.72%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.951'
$ws.Range("E50").Value = '  -2.68%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.13'
$ws.Range("E51").Value = '  +30.32%  '
